$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 246.85075
$ws.Range("I33").Value = 180.81355
$ws.Range("J33").Value = 733.875
$ws.Range("K33").Value = 180.81355
$ws.Range("L33").Value = 733.875
$ws.Range("M33").Value = 48.18645000000001
$ws.Range("N33").Value = -1191.875
$ws.Range("H40").Value = 2718.75
$ws.Range("I40").Value = 1966.6666
$ws.Range("J40").Value = 3685.7144
$ws.Range("K40").Value = 1966.6666
$ws.Range("L40").Value = 3685.7144
$ws.Range("M40").Value = -1791.6666
$ws.Range("N40").Value = -4035.7144
$ws.Range("H43").Value = 2882.353
$ws.Range("I43").Value = 3076.923
$ws.Range("J43").Value = 2250
$ws.Range("K43").Value = 3076.923
$ws.Range("L43").Value = 2250
$ws.Range("M43").Value = -3007.923
$ws.Range("N43").Value = -2388
$ws.Range("H132").Value = 24066.191
$ws.Range("I132").Value = 3760.121
$ws.Range("J132").Value = 98521.78
$ws.Range("K132").Value = 11280.363
$ws.Range("L132").Value = 295565.34
$ws.Range("M132").Value = -8750.363000000001
$ws.Range("N132").Value = -300625.34
$ws.Range("H137").Value = 4090.2979
$ws.Range("I137").Value = 1071.8334
$ws.Range("J137").Value = 4532.0244
$ws.Range("K137").Value = 3215.5002
$ws.Range("L137").Value = 13596.0732
$ws.Range("M137").Value = -665.5001999999999
$ws.Range("N137").Value = -18696.0732

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41656.812
$ws.Range("I32").Value = 43817.027
$ws.Range("K32").Value = 43817.027
$ws.Range("M32").Value = -43530.027
$ws.Range("H110").Value = 1950.3125
$ws.Range("I110").Value = 1875.3572
$ws.Range("J110").Value = 2475
$ws.Range("K110").Value = 1875.3572
$ws.Range("L110").Value = 2475
$ws.Range("M110").Value = 169.6428000000001
$ws.Range("N110").Value = -6565
$ws.Range("H123").Value = 40878
$ws.Range("J123").Value = 40878
$ws.Range("L123").Value = 40878
$ws.Range("N123").Value = -50678
$ws.Range("H132").Value = 14707514
$ws.Range("I132").Value = 20834492
$ws.Range("J132").Value = 2766.5
$ws.Range("K132").Value = 62503476
$ws.Range("L132").Value = 8299.5
$ws.Range("M132").Value = -62500946
$ws.Range("N132").Value = -13359.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3318.58
$ws.Range("I134").Value = 1472.5128
$ws.Range("J134").Value = 4498.8525
$ws.Range("K134").Value = 4417.538399999999
$ws.Range("L134").Value = 13496.5575
$ws.Range("M134").Value = -1882.538399999999
$ws.Range("N134").Value = -18566.5575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J16").Value = 1000
$ws.Range("L16").Value = 1000
$ws.Range("N16").Value = -1574
$ws.Range("J113").Value = 1000
$ws.Range("L113").Value = 1000
$ws.Range("N113").Value = -5340
$ws.Range("H132").Value = 45897.78
$ws.Range("I132").Value = 1733.6818
$ws.Range("J132").Value = 143058.8
$ws.Range("K132").Value = 5201.0454
$ws.Range("L132").Value = 429176.4
$ws.Range("M132").Value = -2671.0454
$ws.Range("N132").Value = -434236.4
$ws.Range("H134").Value = 501662.06
$ws.Range("I134").Value = 1222.3572
$ws.Range("J134").Value = 1002101.8
$ws.Range("K134").Value = 3667.0716
$ws.Range("L134").Value = 3006305.4
$ws.Range("M134").Value = -1132.0716
$ws.Range("N134").Value = -3011375.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 27778036
$ws.Range("I2").Value = 32.2
$ws.Range("J2").Value = 62500540
$ws.Range("K2").Value = 193.2
$ws.Range("L2").Value = 375003240
$ws.Range("M2").Value = -80.20000000000002
$ws.Range("N2").Value = -375003466
$ws.Range("H3").Value = 5450.5713
$ws.Range("I3").Value = 2522.25
$ws.Range("J3").Value = 9355
$ws.Range("K3").Value = 7566.75
$ws.Range("L3").Value = 28065
$ws.Range("M3").Value = -7454.75
$ws.Range("N3").Value = -28289
$ws.Range("H23").Value = 2279.8
$ws.Range("I23").Value = 2000
$ws.Range("J23").Value = 2349.75
$ws.Range("K23").Value = 6000
$ws.Range("L23").Value = 7049.25
$ws.Range("M23").Value = -5765
$ws.Range("N23").Value = -7519.25
$ws.Range("H38").Value = 41833860
$ws.Range("I38").Value = 482
$ws.Range("J38").Value = 71714840
$ws.Range("K38").Value = 1446
$ws.Range("L38").Value = 215144520
$ws.Range("M38").Value = -1099
$ws.Range("N38").Value = -215145214
$ws.Range("H113").Value = 4878.9585
$ws.Range("I113").Value = 8403.923000000001
$ws.Range("J113").Value = 713.0909
$ws.Range("K113").Value = 25211.769
$ws.Range("L113").Value = 2139.2727
$ws.Range("M113").Value = -23041.769
$ws.Range("N113").Value = -6479.2727
$ws.Range("H131").Value = 2187.0454
$ws.Range("J131").Value = 1088.061
$ws.Range("L131").Value = 3264.183
$ws.Range("N131").Value = -13344.183
$ws.Range("H139").Value = 84010.11
$ws.Range("I139").Value = 252780
$ws.Range("J139").Value = 3000.56
$ws.Range("K139").Value = 758340
$ws.Range("L139").Value = 9001.68
$ws.Range("M139").Value = -753200
$ws.Range("N139").Value = -19281.68

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 23937
$ws.Range("J123").Value = 23937
$ws.Range("L123").Value = 23937
$ws.Range("N123").Value = -28837
$ws.Range("H132").Value = 3836.0334
$ws.Range("I132").Value = 1466.8096
$ws.Range("J132").Value = 9364.223
$ws.Range("K132").Value = 4400.4288
$ws.Range("L132").Value = 28092.669
$ws.Range("M132").Value = -1870.4288
$ws.Range("N132").Value = -33152.669

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 39385
$ws.Range("J92").Value = 39385
$ws.Range("L92").Value = 39385
$ws.Range("N92").Value = -44377
$ws.Range("H94").Value = 38220.918
$ws.Range("J94").Value = 38220.918
$ws.Range("L94").Value = 38220.918
$ws.Range("N94").Value = -39572.918
$ws.Range("H132").Value = 3742.8928
$ws.Range("I132").Value = 2739.4443
$ws.Range("J132").Value = 5549.1
$ws.Range("K132").Value = 8218.332900000001
$ws.Range("L132").Value = 16647.3
$ws.Range("M132").Value = -5688.332900000001
$ws.Range("N132").Value = -21707.3

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 35412.5
$ws.Range("J123").Value = 35550
$ws.Range("L123").Value = 35550
$ws.Range("N123").Value = -45350
$ws.Range("H132").Value = 127127.23
$ws.Range("I132").Value = 229122.31
$ws.Range("J132").Value = 2466.5557
$ws.Range("K132").Value = 687366.9299999999
$ws.Range("L132").Value = 7399.6671
$ws.Range("M132").Value = -684836.9299999999
$ws.Range("N132").Value = -12459.6671
